$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This workbook is a rolling "last 10 quarters" income-statement report.
# A new quarter has been added (Q4 ending 1401/12) and the oldest quarter
# (Q2 ending 1399/06) has been dropped, so every data column D..M shifts one
# column to the left, and the newest quarter's data is appended in column M.
# ---------------------------------------------------------------------------

$firstCol = 4   # column D
$lastCol  = 13  # column M

# --- Row 8: quarter-period header labels (shift left, append new quarter) ---
$row8After = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(8, $firstCol + $i).Value = $row8After[$i]
}

# --- Row 9: report-publish-date labels (explicit literal values) ----------
# NOTE: set as literal target values (not a pure shift - one of the middle
# dates was corrected rather than shifted) and guard against Excel
# autoconverting plain "YYYY-MM-DD" text into a date serial, which also
# corrupts the cell style (adds a stray quote-prefix flag). We restore the
# original per-cell formatting afterwards by pasting formats back in from
# the still-intact label cell in column C of the same row.
$row9After = @(
    "1400-10-29 (2)",
    "1401-04-01 (8)",
    "1401-04-28 (2)",
    "1401-08-25 (4)",
    "1401-10-29 (2)",
    "1402-02-30 (8)",
    "1401-04-28",
    "1401-08-25 (2)",
    "1401-10-29",
    "1402-02-30 (2)"
)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(9, $firstCol + $i).Value = "'" + $row9After[$i]
}
$ws.Cells.Item(9, 3).Copy()
$ws.Range("D9:M9").PasteSpecial(-4122)

# --- Rows 11-27: numeric data rows (shift left, append the new quarter) ---
$newLastValues = @{
    11 = 3366125
    12 = -2287793
    13 = 1078332
    14 = -360036
    15 = 0
    16 = 538
    17 = 718834
    18 = -36082
    19 = 6222
    20 = 688974
    21 = 37909
    22 = 726883
    23 = 0
    24 = 726883
    25 = 162
    26 = 4484000
    27 = 162
}

foreach ($row in 11..27) {
    for ($col = $firstCol; $col -lt $lastCol; $col++) {
        $src = $ws.Cells.Item($row, $col + 1)
        $dst = $ws.Cells.Item($row, $col)
        $dst.Value = $src.Value2
    }
    $ws.Cells.Item($row, $lastCol).Value = $newLastValues[$row]
}
